# Update generated output counts (column F) on the "展览" and "全部类型"
# sheets, matching the refreshed scrape numbers committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (1st sheet / sheet1.xml) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1309
$ws1.Range("F5").Value  = 14398
$ws1.Range("F6").Value  = 16788
$ws1.Range("F18").Value = 113
$ws1.Range("F23").Value = 48
$ws1.Range("F24").Value = 31
$ws1.Range("F26").Value = 6870
$ws1.Range("F29").Value = 1141
$ws1.Range("F32").Value = 5782
$ws1.Range("F35").Value = 201
$ws1.Range("F36").Value = 4900
$ws1.Range("F37").Value = 25

# --- Sheet "全部类型" (4th sheet / sheet4.xml) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1309
$ws4.Range("F5").Value  = 14398
$ws4.Range("F6").Value  = 16789
$ws4.Range("F18").Value = 113
$ws4.Range("F24").Value = 48
$ws4.Range("F25").Value = 31
$ws4.Range("F27").Value = 6870
$ws4.Range("F30").Value = 1141
$ws4.Range("F34").Value = 5782
$ws4.Range("F37").Value = 201
$ws4.Range("F38").Value = 4900
$ws4.Range("F39").Value = 25
